$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": PORCELANATO sale for TOSCANO RAMIREZ MONICA CECILIA (row 15)
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M15").Value = -18.25

# Sheet "VENTA MENSUAL": junio (June) sale for same client/row, and the column total
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F15").Value = 1248.94
$wsMensual.Range("F19").Value = 21401.67

# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO group row (16) and TOTAL row (19)
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D16").Value = 20134.48
$wsCumpl.Range("E16").Value = 8075.360000000001
$wsCumpl.Range("F16").Value = 0.7137396029186979

$wsCumpl.Range("D19").Value = 21401.67
$wsCumpl.Range("E19").Value = 25817.63386304603
$wsCumpl.Range("F19").Value = 0.4532398457646261
